$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.886.57'
$ws.Range('E2').Value = '  +1.85%  '
$ws.Range('D3').Value = '1.726.78'
$ws.Range('E3').Value = '  +0.28%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '0.9976'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  -0.25%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '242.00'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -0.47%  '
$ws.Range('E6').Value = '  -0.31%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.4896'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  -0.16%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.2594'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -0.63%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.06212'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +0.29%  '
$ws.Range('D10').Value = '1.729.86'
$ws.Range('E10').Value = '  +0.45%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '16.00'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +3.39%  '
$ws.Range('E12').Value = '  -1.51%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.6091'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +1.69%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '4.485'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -1.68%  '
$ws.Range('E15').Value = '  +0.10%  '
$ws.Range('E16').Value = '  -0.18%  '
$ws.Range('D17').Value = '26.866.57'
$ws.Range('E17').Value = '  +1.78%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.9973'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -0.27%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.000007181'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +0.36%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '11.44'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +0.80%  '
$ws.Range('D21').Value = '1.953.29'
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '4.427'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -0.96%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '8.578'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +0.06%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '5.110'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -0.75%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '138.48'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +0.88%  '
$ws.Range('E26').Value = '  +0.78%  '
$ws.Range('E27').Value = '  +4.68%  '
$ws.Range('E28').Value = '  -1.26%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '106.05'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -0.96%  '
$ws.Range('E30').Value = '  +0.22%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '0.07998'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +0.62%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '3.687'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +0.50%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.04532'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -0.32%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.9969'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -0.28%  '
$ws.Range('E35').Value = '  -0.27%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '1.007'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +1.43%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.6258'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +0.31%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.9366'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +0.56%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '2.057'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +5.96%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '2.451'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +2.37%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.9976'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -0.26%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.01503'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +1.46%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '5.645'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +6.11%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '99.40'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -0.33%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.3853'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +0.48%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '6.892'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +3.02%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.1163'
$ws.Range('D47').ClearFormats()
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.05394'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +0.60%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '7.903'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +2.33%  '
$ws.Range('E50').Value = '  +0.22%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '51.65'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +1.76%  '
